$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("query")

# Insert a new row at row 28 (shifts the existing last data row, "Anna Yukimi
# Yamada", down to row 29) and populate it with the new user's info.
$ws.Rows(28).Insert()

$ws.Range("A28").Value = "Bárbara Port"
$ws.Range("B28").Value = "bport"
$ws.Range("C28").Value = $ws.Range("C29").Value()
$ws.Range("D28").Value = "Corvette, BMW"

# Grow the table (ListObject) to include the new row.
$lo = $ws.ListObjects(1)
$lo.Resize($ws.Range("A1:D29"))

# Update the workbook-level hidden defined name "query" to match new range.
$n = $wb.Names.Item(1)
$n.RefersTo = "=query!`$A`$1:`$D`$29"

# Update selection to match the saved state in the target file.
$ws.Activate()
$ws.Range("E24").Select()
